# issue #5: add legislator_id, name, date into dataframe
# Adds three new columns (date, legislator_name, legislator_id) to the
# "股票" (stocks) worksheet, mirroring the header/value styling already
# used for the existing columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

# New header cells (row 1), styled like the existing header G1.
$ws.Cells.Item(1, 8).Value = "date"
$ws.Cells.Item(1, 9).Value = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"

$ws.Cells.Item(1, 8).Style = $ws.Cells.Item(1, 7).Style
$ws.Cells.Item(1, 9).Style = $ws.Cells.Item(1, 7).Style
$ws.Cells.Item(1, 10).Style = $ws.Cells.Item(1, 7).Style

# New data cells (row 2), styled like the existing data cell G2.
$ws.Cells.Item(2, 8).Value = "2013-12-24"
$ws.Cells.Item(2, 9).Value = "林國正"
$ws.Cells.Item(2, 10).Value = 1742

$ws.Cells.Item(2, 8).Style = $ws.Cells.Item(2, 7).Style
$ws.Cells.Item(2, 9).Style = $ws.Cells.Item(2, 7).Style
$ws.Cells.Item(2, 10).Style = $ws.Cells.Item(2, 7).Style
